$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell C10 value from 18 to 1 (numeric)
$ws.Range("C10").Value = 1
